$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E").Insert()
$ws.Columns("E").ColumnWidth = 32.04
$ws.Range("E4").Value = "Tên  nhóm hàng"
